$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.42115441905784
$ws.Range("D2").Value = 3.617217011492461
$ws.Range("E2").Value = 21.31025839262848
$ws.Range("F2").Value = 26.54696424668288
$ws.Range("G2").Value = 35.33876624063006
$ws.Range("H2").Value = 14.14935312997983
$ws.Range("L2").Value = 9.722047170405057
$ws.Range("N2").Value = 18.35331766152228
$ws.Range("B3").Value = 17.15464674316502
$ws.Range("D3").Value = 3.638905553447395
$ws.Range("E3").Value = 20.96518470046526
$ws.Range("F3").Value = 25.94044271306036
$ws.Range("G3").Value = 34.0595019567597
$ws.Range("H3").Value = 14.04652690696389
$ws.Range("L3").Value = 9.514730164091405
$ws.Range("N3").Value = 18.39138444783763
$ws.Range("B4").Value = 16.99269970616383
$ws.Range("D4").Value = 3.653184789612862
$ws.Range("E4").Value = 20.7493096422299
$ws.Range("F4").Value = 25.57186898515603
$ws.Range("G4").Value = 33.2648307640793
$ws.Range("H4").Value = 13.98785810561196
$ws.Range("L4").Value = 9.386990808304434
$ws.Range("N4").Value = 18.41721839459171
$ws.Range("B5").Value = 16.92721243692192
$ws.Range("D5").Value = 3.659244470587948
$ws.Range("E5").Value = 20.66040194806257
$ws.Range("F5").Value = 25.42288579484907
$ws.Range("G5").Value = 32.93929423936728
$ws.Range("H5").Value = 13.9650929322324
$ws.Range("L5").Value = 9.334896552120913
$ws.Range("N5").Value = 18.42836582764898
$ws.Range("B6").Value = 16.91637137966361
$ws.Range("D6").Value = 3.660265183277748
$ws.Range("E6").Value = 20.64558435290623
$ws.Range("F6").Value = 25.39822789166722
$ws.Range("G6").Value = 32.88515516649034
$ws.Range("H6").Value = 13.96138234281857
$ws.Range("L6").Value = 9.326246120286182
$ws.Range("N6").Value = 18.43025432241594
$ws.Range("B7").Value = 16.99181435855889
$ws.Range("D7").Value = 3.653265539315245
$ws.Range("E7").Value = 20.74811430300517
$ws.Range("F7").Value = 25.56985450032118
$ws.Range("G7").Value = 33.26044651164005
$ws.Range("H7").Value = 13.98754643524373
$ws.Range("L7").Value = 9.386288309193084
$ws.Range("N7").Value = 18.41736622143781
$ws.Range("B8").Value = 17.32896147093782
$ws.Range("D8").Value = 3.624494731803721
$ws.Range("E8").Value = 21.19214734077627
$ws.Range("F8").Value = 26.3371889843387
$ws.Range("G8").Value = 34.89993675576319
$ws.Range("H8").Value = 14.11298401775975
$ws.Range("L8").Value = 9.650699783077265
$ws.Range("N8").Value = 18.36593301912894
$ws.Range("B9").Value = 17.99988663938552
$ws.Range("D9").Value = 3.57576768149265
$ws.Range("E9").Value = 22.02827285113365
$ws.Range("F9").Value = 27.86157884649637
$ws.Range("G9").Value = 38.01762916070555
$ws.Range("H9").Value = 14.39343708729317
$ws.Range("L9").Value = 10.16259257623941
$ws.Range("N9").Value = 18.28454953644862
$ws.Range("B10").Value = 18.49397078989283
$ws.Range("D10").Value = 3.544737182346143
$ws.Range("E10").Value = 22.61776999733273
$ws.Range("F10").Value = 28.97961510577293
$ws.Range("G10").Value = 40.21980668787677
$ws.Range("H10").Value = 14.619051346902
$ws.Range("L10").Value = 10.53071924509348
$ws.Range("N10").Value = 18.23656591456094
$ws.Range("B11").Value = 18.71800379862541
$ws.Range("D11").Value = 3.531674342192471
$ws.Range("E11").Value = 22.87983256441795
$ws.Range("F11").Value = 29.48507131629066
$ws.Range("G11").Value = 41.19730706587404
$ws.Range("H11").Value = 14.72559102555798
$ws.Range("L11").Value = 10.69569538119662
$ws.Range("N11").Value = 18.21728713178598
$ws.Range("B12").Value = 18.80265500337532
$ws.Range("D12").Value = 3.526880815101983
$ws.Range("E12").Value = 22.9781340825802
$ws.Range("F12").Value = 29.67581012071376
$ws.Range("G12").Value = 41.56360608996415
$ws.Range("H12").Value = 14.76646528513221
$ws.Range("L12").Value = 10.75775115469736
$ws.Range("N12").Value = 18.21035215431987
$ws.Range("B13").Value = 18.78443334060206
$ws.Range("D13").Value = 3.527906351412534
$ws.Range("E13").Value = 22.95700561131616
$ws.Range("F13").Value = 29.63476401593389
$ws.Range("G13").Value = 41.48489405980207
$ws.Range("H13").Value = 14.7576392282515
$ws.Range("L13").Value = 10.74440579089444
$ws.Range("N13").Value = 18.21182948862228
$ws.Range("B14").Value = 18.72497226722938
$ws.Range("D14").Value = 3.531276898075168
$ws.Range("E14").Value = 22.88793896908948
$ws.Range("F14").Value = 29.50077804913079
$ws.Range("G14").Value = 41.22752199777454
$ws.Range("H14").Value = 14.72894333966051
$ws.Range("L14").Value = 10.70080943362307
$ws.Range("N14").Value = 18.21670926860233
$ws.Range("B15").Value = 18.68852418061056
$ws.Range("D15").Value = 3.533361441692601
$ws.Range("E15").Value = 22.84551009819628
$ws.Range("F15").Value = 29.41861485024073
$ws.Range("G15").Value = 41.0693614537312
$ws.Range("H15").Value = 14.7114343376612
$ws.Range("L15").Value = 10.6740493858305
$ws.Range("N15").Value = 18.21974583813368
$ws.Range("B16").Value = 18.47930810320779
$ws.Range("D16").Value = 3.545612200659567
$ws.Range("E16").Value = 22.60051573370339
$ws.Range("F16").Value = 28.94650045041221
$ws.Range("G16").Value = 40.15540369573948
$ws.Range("H16").Value = 14.61216465084332
$ws.Range("L16").Value = 10.51988245147211
$ws.Range("N16").Value = 18.2378770318469
$ws.Range("B17").Value = 18.35071610086202
$ws.Range("D17").Value = 3.553398525432224
$ws.Range("E17").Value = 22.44861410225671
$ws.Range("F17").Value = 28.6559114612869
$ws.Range("G17").Value = 39.58822535447778
$ws.Range("H17").Value = 14.5522442752849
$ws.Range("L17").Value = 10.42462556248779
$ws.Range("N17").Value = 18.24965207262473
$ws.Range("B18").Value = 18.27668947410992
$ws.Range("D18").Value = 3.557976032995103
$ws.Range("E18").Value = 22.36067318568414
$ws.Range("F18").Value = 28.48848841083182
$ws.Range("G18").Value = 39.25973292329319
$ws.Range("H18").Value = 14.51814935615872
$ws.Range("L18").Value = 10.36960559166773
$ws.Range("N18").Value = 18.25666477471776
$ws.Range("B19").Value = 18.25161685312981
$ws.Range("D19").Value = 3.559542855575599
$ws.Range("E19").Value = 22.33080164119405
$ws.Range("F19").Value = 28.43175982823179
$ws.Range("G19").Value = 39.14813370235309
$ws.Range("H19").Value = 14.50666984619826
$ws.Range("L19").Value = 10.35093902568889
$ws.Range("N19").Value = 18.25908041010615
$ws.Range("B20").Value = 18.36441214553465
$ws.Range("D20").Value = 3.552559397923392
$ws.Range("E20").Value = 22.46484383101083
$ws.Range("F20").Value = 28.68687613924354
$ws.Range("G20").Value = 39.64883983067777
$ws.Range("H20").Value = 14.55858486582059
$ws.Range("L20").Value = 10.43479014418923
$ws.Range("N20").Value = 18.24837376717323
$ws.Range("B21").Value = 18.74244308003841
$ws.Range("D21").Value = 3.530282718624272
$ws.Range("E21").Value = 22.90825136526601
$ws.Range("F21").Value = 29.54015270423292
$ws.Range("G21").Value = 41.30322585670271
$ws.Range("H21").Value = 14.73735789269243
$ws.Range("L21").Value = 10.71362651583321
$ws.Range("N21").Value = 18.2152660487424
$ws.Range("B22").Value = 18.98839326292216
$ws.Range("D22").Value = 3.516616735575172
$ws.Range("E22").Value = 23.19256390649916
$ws.Range("F22").Value = 30.09384381669164
$ws.Range("G22").Value = 42.36182882161684
$ws.Range("H22").Value = 14.85727012779053
$ws.Range("L22").Value = 10.89340614045731
$ws.Range("N22").Value = 18.19575801225414
$ws.Range("B23").Value = 18.85725285576845
$ws.Range("D23").Value = 3.523828247088102
$ws.Range("E23").Value = 23.04134082658567
$ws.Range("F23").Value = 29.79875935841615
$ws.Range("G23").Value = 41.79901360710332
$ws.Range("H23").Value = 14.79300024502186
$ws.Range("L23").Value = 10.79769775532421
$ws.Range("N23").Value = 18.20597531041334
$ws.Range("B24").Value = 18.35822046089093
$ws.Range("D24").Value = 3.552938452859542
$ws.Range("E24").Value = 22.45750826820518
$ws.Range("F24").Value = 28.67287811190198
$ws.Range("G24").Value = 39.62144351684409
$ws.Range("H24").Value = 14.55571717949215
$ws.Range("L24").Value = 10.43019552811799
$ws.Range("N24").Value = 18.24895093186304
$ws.Range("B25").Value = 17.81783782187338
$ws.Range("D25").Value = 3.588118064661787
$ws.Range("E25").Value = 21.80617733342847
$ws.Range("F25").Value = 27.44855915166089
$ws.Range("G25").Value = 37.1878278147259
$ws.Range("H25").Value = 14.31401873849287
$ws.Range("L25").Value = 10.02523818474811
$ws.Range("N25").Value = 18.30448743018961
